$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 590 of Sheet1 had accidentally picked up a raw date value ("10/10/2018",
# serial 43383, formatted with the "d-mmm" number style) in column A instead
# of the review-text string that belongs there. Fix it the way the author
# did: select the whole row and delete it outright, which shifts every
# following row (the review text in column A / label in column B) up by one
# so the data lines back up correctly. This also removes the now-unused
# trailing row that the shift leaves behind at the end of the sheet.
$ws.Rows.Item(590).Select()
$ws.Rows.Item(590).Delete()
